$wb = $excel.ActiveWorkbook

# Sheet "Alice" holds the start/end request dates that need to become inclusive
# (both start and end now land on the same day).
$alice = $wb.Worksheets.Item("Alice")

$newDate = Get-Date -Year 2019 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0
$alice.Range("A1").Value = $newDate
$alice.Range("B1").Value = $newDate

# Make "Alice" the active/selected sheet (was "Bob" before).
$alice.Activate()
